# Update cryptocurrency price ("Price" column D) values on Sheet1 to reflect
# the latest symbol-list refresh (GitHub Actions run, 2022-12-14 10:46 UTC).
# Values are stored as text in the source sheet, so NumberFormat is forced to
# "@" (Text) before assignment to preserve exact formatting (trailing zeros,
# leading zeros, etc.) instead of letting Excel coerce the string into a
# numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "274.31"
    "D3"  = "22.97"
    "D4"  = "6.343"
    "D5"  = "0.06239"
    "D6"  = "3.658"
    "D7"  = "6.714"
    "D8"  = "1.371"
    "D9"  = "0.8320"
    "D11" = "0.1633"
    "D12" = "0.08291"
    "D13" = "0.03373"
    "D14" = "0.03100"
    "D15" = "0.09310"
    "D16" = "3.884"
    "D17" = "0.001637"
    "D18" = "0.04782"
    "D19" = "0.006378"
    "D20" = "0.005564"
    "D24" = "2.322"
    "D25" = "0.3383"
    "D40" = "0.04696"
    "D41" = "0.007031"
    "D43" = "0.003600"
    "D45" = "0.00006253"
    "D47" = "0.9000"
    "D48" = "0.03023"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
